# Update cryptos list: refresh Price (column D) and Volume(1h) (column E)
# values for rows 2-51 on the active worksheet.
#
# Some Price values are numeric-looking strings (e.g. "155.56", "0.1000")
# that must stay as literal text (preserving formats such as trailing
# zeros). Setting NumberFormat to "@" (Text) before assigning the Value
# forces Excel to store them as text instead of auto-converting them to
# numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.528.13'
$ws.Range("E2").Value = '  -0.56%  '
$ws.Range("D3").Value = '2.646.46'
$ws.Range("E3").Value = '  -1.38%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("E5").Value = '  -1.47%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '155.56'
$ws.Range("E6").Value = '  -0.81%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  +5.66%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.126'
$ws.Range("E9").Value = '  +2.59%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.394'
$ws.Range("E10").Value = '  -1.18%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.77'
$ws.Range("E11").Value = '  -3.18%  '
$ws.Range("E12").Value = '  +0.16%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '28.67'
$ws.Range("E13").Value = '  -2.97%  '
$ws.Range("E14").Value = '  -3.12%  '
$ws.Range("D15").Value = '3.121.08'
$ws.Range("E15").Value = '  -1.42%  '
$ws.Range("D16").Value = '65.393.52'
$ws.Range("E16").Value = '  -0.45%  '
$ws.Range("D17").Value = '2.629.41'
$ws.Range("E17").Value = '  -1.65%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.59'
$ws.Range("E18").Value = '  -0.30%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.73'
$ws.Range("E19").Value = '  -2.26%  '
$ws.Range("E20").Value = '  -2.14%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '347.96'
$ws.Range("E21").Value = '  -1.39%  '
$ws.Range("E22").Value = '  -0.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.94'
$ws.Range("E23").Value = '  -2.26%  '
$ws.Range("E24").Value = '  +0.24%  '
$ws.Range("E25").Value = '  -2.81%  '
$ws.Range("E26").Value = '  +0.34%  '
$ws.Range("E27").Value = '  -3.21%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.164'
$ws.Range("E28").Value = '  -3.45%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.31%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '538.53'
$ws.Range("E30").Value = '  +0.58%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.88'
$ws.Range("E31").Value = '  -3.70%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.12'
$ws.Range("E32").Value = '  -2.76%  '
$ws.Range("E33").Value = '  -1.53%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.37'
$ws.Range("E34").Value = '  -3.39%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.39'
$ws.Range("E35").Value = '  -0.47%  '
$ws.Range("E36").Value = '  -1.84%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '20.25'
$ws.Range("E37").Value = '  -1.25%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  +0.05%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '154.84'
$ws.Range("E39").Value = '  -3.56%  '
$ws.Range("E40").Value = '  -3.22%  '
$ws.Range("E41").Value = '  -0.04%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '160.25'
$ws.Range("E42").Value = '  -3.88%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.06'
$ws.Range("E43").Value = '  -1.39%  '
$ws.Range("E44").Value = '  +2.67%  '
$ws.Range("E45").Value = '  -3.06%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.48'
$ws.Range("E46").Value = '  -3.16%  '
$ws.Range("E47").Value = '  -2.79%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0253'
$ws.Range("E48").Value = '  -4.16%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.1000'
$ws.Range("E49").Value = '  +0.62%  '
$ws.Range("E50").Value = '  +7.07%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.55'
$ws.Range("E51").Value = '  -3.94%  '
